$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"1.174437034476609"
$ws.Range("C2").Value = [double]"-3.820608917710134"
$ws.Range("D2").Value = [double]"0.01383533758484761"
$ws.Range("E2").Value = [double]"0.3121335553179975"
$ws.Range("F2").Value = [double]"0"
$ws.Range("G2").Value = [double]"2.883958458594906e-32"
$ws.Range("H2").Value = [double]"0.874617094785566"
$ws.Range("I2").Value = [double]"3.446465529938971"
$ws.Range("J2").Value = [double]"-1.407355177109377e-14"
$ws.Range("K2").Value = [double]"1035"
$ws.Range("B3").Value = [double]"1.11444761042579"
$ws.Range("C3").Value = [double]"-2.989473328999499"
$ws.Range("D3").Value = [double]"0.01568101456106926"
$ws.Range("E3").Value = [double]"0.3490895151673853"
$ws.Range("F3").Value = [double]"0"
$ws.Range("G3").Value = [double]"3.9410857620688e-17"
$ws.Range("H3").Value = [double]"0.8307545166893965"
$ws.Range("I3").Value = [double]"3.524729130825958"
$ws.Range("J3").Value = [double]"1.168846246216421e-14"
$ws.Range("K3").Value = [double]"1031"
$ws.Range("B4").Value = [double]"1.099797084518694"
$ws.Range("C4").Value = [double]"-3.178141876699287"
$ws.Range("D4").Value = [double]"0.0183317341364989"
$ws.Range("E4").Value = [double]"0.4291734904402102"
$ws.Range("F4").Value = [double]"0"
$ws.Range("G4").Value = [double]"2.708797722942366e-13"
$ws.Range("H4").Value = [double]"0.7766652926778326"
$ws.Range("I4").Value = [double]"4.330301978577632"
$ws.Range("J4").Value = [double]"2.625650687977535e-14"
$ws.Range("K4").Value = [double]"1037"
$ws.Range("B5").Value = [double]"1.112843991707752"
$ws.Range("C5").Value = [double]"-4.090832345812087"
$ws.Range("D5").Value = [double]"0.01495658566571337"
$ws.Range("E5").Value = [double]"0.3676741405347682"
$ws.Range("F5").Value = [double]"0"
$ws.Range("G5").Value = [double]"2.631266171151702e-27"
$ws.Range("H5").Value = [double]"0.8361298594968598"
$ws.Range("I5").Value = [double]"3.523094826418407"
$ws.Range("J5").Value = [double]"9.138350916215455e-15"
$ws.Range("K5").Value = [double]"1087"
$ws.Range("B6").Value = [double]"1.176558507630624"
$ws.Range("C6").Value = [double]"-3.976824221892961"
$ws.Range("D6").Value = [double]"0.01703825149145682"
$ws.Range("E6").Value = [double]"0.3962774028778818"
$ws.Range("F6").Value = [double]"0"
$ws.Range("G6").Value = [double]"1.023991325948511e-22"
$ws.Range("H6").Value = [double]"0.8158936334077829"
$ws.Range("I6").Value = [double]"3.907314871203506"
$ws.Range("J6").Value = [double]"1.705829870266363e-14"
$ws.Range("K6").Value = [double]"1078"
$ws.Range("B7").Value = [double]"1.219614365170474"
$ws.Range("C7").Value = [double]"-4.656912389355643"
$ws.Range("D7").Value = [double]"0.02031227281312861"
$ws.Range("E7").Value = [double]"0.4656137010230384"
$ws.Range("F7").Value = [double]"0"
$ws.Range("G7").Value = [double]"1.66041635863134e-22"
$ws.Range("H7").Value = [double]"0.7837045324610361"
$ws.Range("I7").Value = [double]"4.419506898969251"
$ws.Range("J7").Value = [double]"4.2760846685663e-16"
$ws.Range("K7").Value = [double]"997"
$ws.Range("B8").Value = [double]"1.068442746134074"
$ws.Range("C8").Value = [double]"-0.6424991007790973"
$ws.Range("D8").Value = [double]"0.01507295218724373"
$ws.Range("E8").Value = [double]"0.334342504443856"
$ws.Range("F8").Value = [double]"0"
$ws.Range("G8").Value = [double]"0.05500647903070693"
$ws.Range("H8").Value = [double]"0.8644336197950682"
$ws.Range("I8").Value = [double]"3.223113065141798"
$ws.Range("J8").Value = [double]"1.412091259675136e-15"
$ws.Range("K8").Value = [double]"790"
$ws.Range("B9").Value = [double]"1.107730834345226"
$ws.Range("C9").Value = [double]"-2.764929522411111"
$ws.Range("D9").Value = [double]"0.01961420507900099"
$ws.Range("E9").Value = [double]"0.4384179208384364"
$ws.Range("F9").Value = [double]"1.578938943504622e-317"
$ws.Range("G9").Value = [double]"4.22726191307851e-10"
$ws.Range("H9").Value = [double]"0.7560763277210523"
$ws.Range("I9").Value = [double]"4.857762171809617"
$ws.Range("J9").Value = [double]"3.302541988130359e-14"
$ws.Range("K9").Value = [double]"1031"
